# Applies the weekly Fruta/Hortaliza data refresh for the Espárragos
# subset (Mercado Mayorista Lo Valledor de Santiago): updated Fecha,
# Volumen/Precio figures, Unidad de comercializacion and Origen for
# rows 2-8, 12-21 (rows 9-11 are unchanged in this week's extract).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44162
$ws.Range("J2").Value = 5200
$ws.Range("K2").Value = 1100
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = 1100
$ws.Range("O2").Value = 'Provincia de Linares'
$ws.Range("P2").Value = 1100

# Row 3
$ws.Range("D3").Value = 44162
$ws.Range("J3").Value = 3400
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = 900
$ws.Range("O3").Value = 'Provincia de Linares'
$ws.Range("P3").Value = 900

# Row 4
$ws.Range("D4").Value = 44160
$ws.Range("H4").Value = 'Verde'
$ws.Range("J4").Value = 210
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 13000
$ws.Range("N4").Value = '$/bandeja 10 kilos'
$ws.Range("O4").Value = 'Región Metropolitana'
$ws.Range("P4").Value = 1300

# Row 5
$ws.Range("D5").Value = 44160
$ws.Range("H5").Value = 'Verde'
$ws.Range("J5").Value = 340
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("N5").Value = '$/bandeja 10 kilos'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 1100

# Row 6
$ws.Range("D6").Value = 44160
$ws.Range("H6").Value = 'Verde'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 4300
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 1200
$ws.Range("M6").Value = 1200
$ws.Range("N6").Value = '$/kilo'
$ws.Range("O6").Value = 'Región Metropolitana'
$ws.Range("P6").Value = 1200
$ws.Range("Q6").Value = 1

# Row 7
$ws.Range("D7").Value = 44160
$ws.Range("I7").Value = 'Segunda'
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 9000
$ws.Range("N7").Value = '$/bandeja 10 kilos'
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("P7").Value = 900
$ws.Range("Q7").Value = 10

# Row 8
$ws.Range("D8").Value = 44160
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 1000
$ws.Range("O8").Value = 'Región Metropolitana'
$ws.Range("P8").Value = 1000

# Row 12
$ws.Range("D12").Value = 44467
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 2000
$ws.Range("N12").Value = '$/caja 10 kilos'
$ws.Range("P12").Value = 200
$ws.Range("Q12").Value = 10

# Row 13
$ws.Range("D13").Value = 44467
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1500
$ws.Range("N13").Value = '$/caja 10 kilos'
$ws.Range("P13").Value = 150
$ws.Range("Q13").Value = 10

# Row 14
$ws.Range("D14").Value = 44467
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 1300
$ws.Range("L14").Value = 1300
$ws.Range("M14").Value = 1300
$ws.Range("N14").Value = '$/caja 10 kilos'
$ws.Range("P14").Value = 130
$ws.Range("Q14").Value = 10

# Row 15
$ws.Range("D15").Value = 44161
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 4300
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 1000
$ws.Range("N15").Value = '$/kilo'
$ws.Range("O15").Value = 'Provincia de Linares'
$ws.Range("P15").Value = 1000
$ws.Range("Q15").Value = 1

# Row 16
$ws.Range("D16").Value = 44161
$ws.Range("I16").Value = 'Segunda'
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = 800
$ws.Range("N16").Value = '$/kilo'
$ws.Range("O16").Value = 'Provincia de Linares'
$ws.Range("P16").Value = 800
$ws.Range("Q16").Value = 1

# Row 17
$ws.Range("D17").Value = 44474
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("I17").Value = 'Banquete'
$ws.Range("J17").Value = 780
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 1600
$ws.Range("M17").Value = 1558
$ws.Range("O17").Value = 'Provincia de Linares'
$ws.Range("P17").Value = 1558

# Row 18
$ws.Range("D18").Value = 44474
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 520
$ws.Range("K18").Value = 1300
$ws.Range("L18").Value = 1400
$ws.Range("M18").Value = 1348
$ws.Range("N18").Value = '$/kilo'
$ws.Range("O18").Value = 'Provincia de Linares'
$ws.Range("P18").Value = 1348
$ws.Range("Q18").Value = 1

# Row 19
$ws.Range("D19").Value = 44474
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("J19").Value = 400
$ws.Range("L19").Value = 1200
$ws.Range("M19").Value = 1100
$ws.Range("O19").Value = 'Provincia de Linares'
$ws.Range("P19").Value = 1100

# Row 20
$ws.Range("D20").Value = 44159
$ws.Range("J20").Value = 4300
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = 1000
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("P20").Value = 1000

# Row 21
$ws.Range("D21").Value = 44159
$ws.Range("J21").Value = 2500
$ws.Range("K21").Value = 800
$ws.Range("L21").Value = 800
$ws.Range("M21").Value = 800
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 800
